$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.081008657837123
$ws.Cells.Item(2, 4).Value = 1.085986822212265
$ws.Cells.Item(2, 5).Value = 1.084533756139163
$ws.Cells.Item(2, 6).Value = 1.09615399901796
$ws.Cells.Item(2, 9).Value = 1.050697501367217
$ws.Cells.Item(2, 10).Value = 1.085885365653844
$ws.Cells.Item(2, 11).Value = 1.088645832764586
$ws.Cells.Item(2, 12).Value = 1.08719652855668
$ws.Cells.Item(2, 13).Value = 1.098786988442533
$ws.Cells.Item(2, 14).Value = 1.087427447425865
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.082593907448944
$ws.Cells.Item(3, 4).Value = 1.087441188987589
$ws.Cells.Item(3, 5).Value = 1.085961616526383
$ws.Cells.Item(3, 6).Value = 1.097702308494605
$ws.Cells.Item(3, 9).Value = 1.051069649742426
$ws.Cells.Item(3, 10).Value = 1.087129235136191
$ws.Cells.Item(3, 11).Value = 1.089917968186072
$ws.Cells.Item(3, 12).Value = 1.0884419465658
$ws.Cells.Item(3, 13).Value = 1.10015474460595
$ws.Cells.Item(3, 14).Value = 1.088673083345551
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.083617941740146
$ws.Cells.Item(4, 4).Value = 1.088380876010686
$ws.Cells.Item(4, 5).Value = 1.086884220969942
$ws.Cells.Item(4, 6).Value = 1.098702973102083
$ws.Cells.Item(4, 9).Value = 1.051307837249775
$ws.Cells.Item(4, 10).Value = 1.08793194349308
$ws.Cells.Item(4, 11).Value = 1.090739213529552
$ws.Cells.Item(4, 12).Value = 1.089245967842914
$ws.Cells.Item(4, 13).Value = 1.101038067709939
$ws.Cells.Item(4, 14).Value = 1.089476931640378
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.084048040644832
$ws.Cells.Item(5, 4).Value = 1.088775595929261
$ws.Cells.Item(5, 5).Value = 1.087271775663825
$ws.Cells.Item(5, 6).Value = 1.099123373991352
$ws.Cells.Item(5, 9).Value = 1.05140734703139
$ws.Cells.Item(5, 10).Value = 1.088268891980627
$ws.Cells.Item(5, 11).Value = 1.091084014139931
$ws.Cells.Item(5, 12).Value = 1.08958354259694
$ws.Cells.Item(5, 13).Value = 1.10140901518562
$ws.Cells.Item(5, 14).Value = 1.089814358633426
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.084120232637107
$ws.Cells.Item(6, 4).Value = 1.088841852319349
$ws.Cells.Item(6, 5).Value = 1.087336829931167
$ws.Cells.Item(6, 6).Value = 1.099193945038669
$ws.Cells.Item(6, 9).Value = 1.05142401862352
$ws.Cells.Item(6, 10).Value = 1.088325437392734
$ws.Cells.Item(6, 11).Value = 1.091141881428989
$ws.Cells.Item(6, 12).Value = 1.089640197516387
$ws.Cells.Item(6, 13).Value = 1.101471275554155
$ws.Cells.Item(6, 14).Value = 1.089870984346505
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.083623690323932
$ws.Cells.Item(7, 4).Value = 1.088386151543772
$ws.Cells.Item(7, 5).Value = 1.086889400697156
$ws.Cells.Item(7, 6).Value = 1.098708591604309
$ws.Cells.Item(7, 9).Value = 1.05130916935433
$ws.Cells.Item(7, 10).Value = 1.087936447811492
$ws.Cells.Item(7, 11).Value = 1.090743822535853
$ws.Cells.Item(7, 12).Value = 1.089250480237661
$ws.Cells.Item(7, 13).Value = 1.101043025897477
$ws.Cells.Item(7, 14).Value = 1.089481442355439
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.081544763411392
$ws.Cells.Item(8, 4).Value = 1.08647862350074
$ws.Cells.Item(8, 5).Value = 1.085016585106685
$ws.Cells.Item(8, 6).Value = 1.096677509515869
$ws.Cells.Item(8, 9).Value = 1.050823814025494
$ws.Cells.Item(8, 10).Value = 1.086306188069514
$ws.Cells.Item(8, 11).Value = 1.089076156656211
$ws.Cells.Item(8, 12).Value = 1.087617809852877
$ws.Cells.Item(8, 13).Value = 1.099249585633407
$ws.Cells.Item(8, 14).Value = 1.087848867457634
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.077867804430607
$ws.Cells.Item(9, 4).Value = 1.083106380026237
$ws.Cells.Item(9, 5).Value = 1.081706044906722
$ws.Cells.Item(9, 6).Value = 1.093088987341021
$ws.Cells.Item(9, 9).Value = 1.049948403133732
$ws.Cells.Item(9, 10).Value = 1.083416633633492
$ws.Cells.Item(9, 11).Value = 1.086122578427304
$ws.Cells.Item(9, 12).Value = 1.084726394391101
$ws.Cells.Item(9, 13).Value = 1.096075928924237
$ws.Cells.Item(9, 14).Value = 1.084955209522876
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.075406792884788
$ws.Cells.Item(10, 4).Value = 1.080850409747116
$ws.Cells.Item(10, 5).Value = 1.079491581990393
$ws.Cells.Item(10, 6).Value = 1.090689772807788
$ws.Cells.Item(10, 9).Value = 1.049351094698906
$ws.Cells.Item(10, 10).Value = 1.081478529996803
$ws.Cells.Item(10, 11).Value = 1.084143061634427
$ws.Cells.Item(10, 12).Value = 1.082788662914622
$ws.Cells.Item(10, 13).Value = 1.093950714173298
$ws.Cells.Item(10, 14).Value = 1.08301435355672
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.074338723513763
$ws.Cells.Item(11, 4).Value = 1.079871595539167
$ws.Cells.Item(11, 5).Value = 1.07853083178169
$ws.Cells.Item(11, 6).Value = 1.089649149765461
$ws.Cells.Item(11, 9).Value = 1.049089169413289
$ws.Cells.Item(11, 10).Value = 1.08063643236617
$ws.Cells.Item(11, 11).Value = 1.083283332388358
$ws.Cells.Item(11, 12).Value = 1.081947110862081
$ws.Cells.Item(11, 13).Value = 1.093028133267266
$ws.Cells.Item(11, 14).Value = 1.082171060050861
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.07394161852334
$ws.Cells.Item(12, 4).Value = 1.07950771636147
$ws.Cells.Item(12, 5).Value = 1.078173676286624
$ws.Cells.Item(12, 6).Value = 1.089262344118987
$ws.Cells.Item(12, 9).Value = 1.048991381843327
$ws.Cells.Item(12, 10).Value = 1.080323198007798
$ws.Cells.Item(12, 11).Value = 1.082963593995038
$ws.Cells.Item(12, 12).Value = 1.081634137084223
$ws.Cells.Item(12, 13).Value = 1.092685084072512
$ws.Cells.Item(12, 14).Value = 1.081857380863772
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.07402681609935
$ws.Cells.Item(13, 4).Value = 1.079585783580681
$ws.Cells.Item(13, 5).Value = 1.078250300619926
$ws.Cells.Item(13, 6).Value = 1.08934532769973
$ws.Cells.Item(13, 9).Value = 1.04901238014224
$ws.Cells.Item(13, 10).Value = 1.080390407966315
$ws.Cells.Item(13, 11).Value = 1.083032197049159
$ws.Cells.Item(13, 12).Value = 1.08170128852314
$ws.Cells.Item(13, 13).Value = 1.092758685791109
$ws.Cells.Item(13, 14).Value = 1.081924686268139
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.074305906433014
$ws.Cells.Item(14, 4).Value = 1.079841523426093
$ws.Cells.Item(14, 5).Value = 1.07850131516523
$ws.Cells.Item(14, 6).Value = 1.08961718189171
$ws.Cells.Item(14, 9).Value = 1.049081096415766
$ws.Cells.Item(14, 10).Value = 1.08061054937688
$ws.Cells.Item(14, 11).Value = 1.083256910852956
$ws.Cells.Item(14, 12).Value = 1.081921248229722
$ws.Cells.Item(14, 13).Value = 1.092999784140145
$ws.Cells.Item(14, 14).Value = 1.082145140304758
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.074477812993498
$ws.Cells.Item(15, 4).Value = 1.079999052634826
$ws.Cells.Item(15, 5).Value = 1.078655934810668
$ws.Cells.Item(15, 6).Value = 1.089784643917726
$ws.Cells.Item(15, 9).Value = 1.049123368849528
$ws.Cells.Item(15, 10).Value = 1.08074612697499
$ws.Cells.Item(15, 11).Value = 1.083395311633225
$ws.Cells.Item(15, 12).Value = 1.082056721555324
$ws.Cells.Item(15, 13).Value = 1.093148284618123
$ws.Cells.Item(15, 14).Value = 1.082280910438608
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.075477624921626
$ws.Cells.Item(16, 4).Value = 1.080915328281225
$ws.Cells.Item(16, 5).Value = 1.07955530361413
$ws.Cells.Item(16, 6).Value = 1.090758797917474
$ws.Cells.Item(16, 9).Value = 1.049368408308746
$ws.Cells.Item(16, 10).Value = 1.081534355739272
$ws.Cells.Item(16, 11).Value = 1.084200063859495
$ws.Cells.Item(16, 12).Value = 1.08284446056298
$ws.Cells.Item(16, 13).Value = 1.094011892602515
$ws.Cells.Item(16, 14).Value = 1.083070258578147
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.076104121107412
$ws.Cells.Item(17, 4).Value = 1.081489551657242
$ws.Cells.Item(17, 5).Value = 1.08011894631023
$ws.Cells.Item(17, 6).Value = 1.091369384567647
$ws.Cells.Item(17, 9).Value = 1.049521232975542
$ws.Cells.Item(17, 10).Value = 1.082028012383508
$ws.Cells.Item(17, 11).Value = 1.084704165685229
$ws.Cells.Item(17, 12).Value = 1.083337913039003
$ws.Cells.Item(17, 13).Value = 1.094552976109807
$ws.Cells.Item(17, 14).Value = 1.08356461627144
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.076469311500813
$ws.Cells.Item(18, 4).Value = 1.081824297527541
$ws.Cells.Item(18, 5).Value = 1.08044752941692
$ws.Cells.Item(18, 6).Value = 1.091725361479467
$ws.Cells.Item(18, 9).Value = 1.04961005609887
$ws.Cells.Item(18, 10).Value = 1.082315676057917
$ws.Cells.Item(18, 11).Value = 1.084997950857878
$ws.Cells.Item(18, 12).Value = 1.083625494705611
$ws.Cells.Item(18, 13).Value = 1.094868354921377
$ws.Cells.Item(18, 14).Value = 1.083852688461263
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.076593792577355
$ws.Cells.Item(19, 4).Value = 1.08193840544415
$ws.Cells.Item(19, 5).Value = 1.080559537470634
$ws.Cells.Item(19, 6).Value = 1.091846712202595
$ws.Cells.Item(19, 9).Value = 1.0496402888252
$ws.Cells.Item(19, 10).Value = 1.082413715048602
$ws.Cells.Item(19, 11).Value = 1.085098082040991
$ws.Cells.Item(19, 12).Value = 1.08372351201804
$ws.Cells.Item(19, 13).Value = 1.094975852866016
$ws.Cells.Item(19, 14).Value = 1.08395086667856
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.076036928300403
$ws.Cells.Item(20, 4).Value = 1.081427962527694
$ws.Cells.Item(20, 5).Value = 1.080058491478916
$ws.Cells.Item(20, 6).Value = 1.091303891838651
$ws.Cells.Item(20, 9).Value = 1.049504869137664
$ws.Cells.Item(20, 10).Value = 1.081975076500922
$ws.Cells.Item(20, 11).Value = 1.084650106134024
$ws.Cells.Item(20, 12).Value = 1.083284995225131
$ws.Cells.Item(20, 13).Value = 1.094494946413826
$ws.Cells.Item(20, 14).Value = 1.083511605213829
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.07422373176406
$ws.Cells.Item(21, 4).Value = 1.07976622287952
$ws.Cells.Item(21, 5).Value = 1.078427405715932
$ws.Cells.Item(21, 6).Value = 1.089537135190045
$ws.Cells.Item(21, 9).Value = 1.049060874926293
$ws.Cells.Item(21, 10).Value = 1.080545735480694
$ws.Cells.Item(21, 11).Value = 1.083190749255237
$ws.Cells.Item(21, 12).Value = 1.081856486240316
$ws.Cells.Item(21, 13).Value = 1.092928796736671
$ws.Cells.Item(21, 14).Value = 1.082080234365405
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.073081519668095
$ws.Cells.Item(22, 4).Value = 1.078719657895916
$ws.Cells.Item(22, 5).Value = 1.077400194636566
$ws.Cells.Item(22, 6).Value = 1.08842472723106
$ws.Cells.Item(22, 9).Value = 1.048778841946012
$ws.Cells.Item(22, 10).Value = 1.079644491782671
$ws.Cells.Item(22, 11).Value = 1.082270894505708
$ws.Cells.Item(22, 12).Value = 1.080956100920828
$ws.Cells.Item(22, 13).Value = 1.091941999998997
$ws.Cells.Item(22, 14).Value = 1.081177710797956
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.073687238447824
$ws.Cells.Item(23, 4).Value = 1.079274631929036
$ws.Cells.Item(23, 5).Value = 1.077944901081823
$ws.Cells.Item(23, 6).Value = 1.089014588466464
$ws.Cells.Item(23, 9).Value = 1.048928626623687
$ws.Cells.Item(23, 10).Value = 1.080122503566425
$ws.Cells.Item(23, 11).Value = 1.082758747616988
$ws.Cells.Item(23, 12).Value = 1.081433625825983
$ws.Cells.Item(23, 13).Value = 1.092465321238345
$ws.Cells.Item(23, 14).Value = 1.081656401413273
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.07606729056155
$ws.Cells.Item(24, 4).Value = 1.0814557925882
$ws.Cells.Item(24, 5).Value = 1.080085808971783
$ws.Cells.Item(24, 6).Value = 1.091333485702965
$ws.Cells.Item(24, 9).Value = 1.049512264230984
$ws.Cells.Item(24, 10).Value = 1.081998996806521
$ws.Cells.Item(24, 11).Value = 1.084674534087499
$ws.Cells.Item(24, 12).Value = 1.083308907251682
$ws.Cells.Item(24, 13).Value = 1.094521168234241
$ws.Cells.Item(24, 14).Value = 1.083535559489006
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.078820056872267
$ws.Cells.Item(25, 4).Value = 1.083979528361682
$ws.Cells.Item(25, 5).Value = 1.082563176800362
$ws.Cells.Item(25, 6).Value = 1.094017881024008
$ws.Cells.Item(25, 9).Value = 1.050177121531633
$ws.Cells.Item(25, 10).Value = 1.084165689934158
$ws.Cells.Item(25, 11).Value = 1.086887962554953
$ws.Cells.Item(25, 12).Value = 1.085475648143184
$ws.Cells.Item(25, 13).Value = 1.096898026994317
$ws.Cells.Item(25, 14).Value = 1.085705329569408
